# Apply "Horarios actualizados" update to every sheet in the workbook.
# For each sheet:
#   - Update the "Última actualización: HH:MM:SS" cell (row 2) to the new timestamp
#   - Increment the "Total filas: N" counter (row 3) by 1
#   - Append one new data row with the new scrape values

$wb = $excel.ActiveWorkbook

$newTimestamp = "02:47:42"
$newArrival   = "04:45"
$newLinea     = "215A_EL PATO"
$newMinutos   = 118

foreach ($ws in $wb.Worksheets) {

    # --- Update "Última actualización" (row 2, column A) ---
    $ws.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"

    # --- Update "Total filas" (row 3, column A) ---
    $totalFilasCell = $ws.Cells.Item(3, 1)
    $currentText = [string]$totalFilasCell.Value()
    $currentCount = [int]($currentText -replace '[^0-9]', '')
    $newCount = $currentCount + 1
    $totalFilasCell.Value = "Total filas: $newCount"

    # --- Find last used data row (header row is row 5, data starts row 6) ---
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
    $newRow = $lastRow + 1

    # --- Append the new data row ---
    $ws.Cells.Item($newRow, 1).Value = $newTimestamp
    $ws.Cells.Item($newRow, 2).Value = $newArrival
    $ws.Cells.Item($newRow, 3).Value = $newLinea
    $ws.Cells.Item($newRow, 4).Value = $newMinutos
}
